$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = '30.866.68'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.915.94'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = $defaultStyle
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'239.65"
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  -3.33%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = $defaultStyle
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').Value = "'0.4906"
$ws.Range('D7').Style = $defaultStyle
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('D8').Value = "'0.2967"
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.06764"
$ws.Range('D9').Style = $defaultStyle
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').Value = '1.876.50'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').Value = "'17.01"
$ws.Range('D11').Style = $defaultStyle
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').Value = "'0.07313"
$ws.Range('D12').Style = $defaultStyle
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = "'5.156"
$ws.Range('D13').Style = $defaultStyle
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').Value = "'89.96"
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = "'0.6737"
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D16').Value = '30.814.25'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = "'0.000007939"
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = "'13.49"
$ws.Range('D18').Style = $defaultStyle
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '2.143.22'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').Value = "'5.172"
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').Value = '  +6.65%  '
$ws.Range('D23').Value = "'205.45"
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').Value = '  +7.36%  '
$ws.Range('D24').Value = "'6.253"
$ws.Range('D24').Style = $defaultStyle
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('D25').Value = "'9.681"
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').Value = '  +2.84%  '
$ws.Range('D26').Value = "'158.83"
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  +1.41%  '
$ws.Range('D27').Value = "'18.89"
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('D28').Value = "'1.975"
$ws.Range('D28').Style = $defaultStyle
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').Value = "'4.334"
$ws.Range('D30').Style = $defaultStyle
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = "'0.09195"
$ws.Range('D31').Style = $defaultStyle
$ws.Range('E31').Value = '  +2.16%  '
$ws.Range('D32').Value = "'4.061"
$ws.Range('D32').Style = $defaultStyle
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').Value = "'0.05191"
$ws.Range('D33').Style = $defaultStyle
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').Value = "'0.7523"
$ws.Range('D34').Style = $defaultStyle
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('D35').Value = "'1.118"
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('D36').Value = "'2.718"
$ws.Range('D36').Style = $defaultStyle
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = "'0.01858"
$ws.Range('D37').Style = $defaultStyle
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').Value = "'2.730"
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Value = '  +1.58%  '
$ws.Range('D39').Value = "'0.9257"
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('D40').Value = "'2.090"
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('D41').Value = "'0.4506"
$ws.Range('D41').Style = $defaultStyle
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('D42').Value = "'107.75"
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Value = '  +1.42%  '
$ws.Range('D43').Value = "'5.892"
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Value = '  +2.17%  '
$ws.Range('D44').Value = "'1.008"
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').Value = "'0.1394"
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').Value = '  +3.61%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'69.01"
$ws.Range('D46').Style = $defaultStyle
$ws.Range('E46').Value = '  +18.73%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = "'7.717"
$ws.Range('D47').Style = $defaultStyle
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = "'35.84"
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').Value = '  +6.58%  '
$ws.Range('D49').Value = "'0.4097"
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').Value = '  +3.53%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'8.974"
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').Value = '  +2.99%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.05942"
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Value = '  +1.34%  '
